{"js": "// Update the answer text in the \"two-digit divided by one-digit\" worksheet\n// table. The worksheet table has 5 \"data\" rows (0, 4, 8, 12, 16 - 0 based)\n// each holding 5 division-problem cells; the remaining rows are blank\n// answer rows. Every data cell's text is replaced with a new value,\n// matching the cell's (row, column) position so the edit is unambiguous\n// even where an old/new value collides with another cell's value.\n\nconst rowReplacements = [\n  { row: 0, values: [\"99\u00f74=24, 3\", \"26\u00f76=4, 2\", \"42\u00f73=14, 0\", \"73\u00f77=10, 3\", \"76\u00f78=9, 4\"] },\n  { row: 4, values: [\"61\u00f76=10, 1\", \"46\u00f74=11, 2\", \"34\u00f76=5, 4\", \"47\u00f74=11, 3\", \"75\u00f73=25, 0\"] },\n  { row: 8, values: [\"48\u00f78=6, 0\", \"24\u00f72=12, 0\", \"99\u00f78=12, 3\", \"43\u00f74=10, 3\", \"67\u00f75=13, 2\"] },\n  { row: 12, values: [\"69\u00f75=13, 4\", \"37\u00f73=12, 1\", \"82\u00f78=10, 2\", \"60\u00f76=10, 0\", \"97\u00f74=24, 1\"] },\n  { row: 16, values: [\"51\u00f77=7, 2\", \"84\u00f74=21, 0\", \"18\u00f79=2, 0\", \"39\u00f79=4, 3\", \"69\u00f78=8, 5\"] },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const { row, values } of rowReplacements) {\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    const paragraph = cell.body.paragraphs.getFirst();\n    paragraph.insertText(values[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the answer text in the \"two-digit divided by one-digit\" worksheet\n# table. The worksheet table has 5 \"data\" rows (Word 1-based rows 1, 5, 9,\n# 13, 17) each holding 5 division-problem cells; the rows in between are\n# blank answer rows. Every data cell's text is replaced with a new value,\n# addressed by its exact (row, column) position so the edit is unambiguous\n# even where an old/new value collides with another cell's value.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowReplacements = @(\n    @{ Row = 1;  Values = @(\"99\u00f74=24, 3\", \"26\u00f76=4, 2\", \"42\u00f73=14, 0\", \"73\u00f77=10, 3\", \"76\u00f78=9, 4\") },\n    @{ Row = 5;  Values = @(\"61\u00f76=10, 1\", \"46\u00f74=11, 2\", \"34\u00f76=5, 4\", \"47\u00f74=11, 3\", \"75\u00f73=25, 0\") },\n    @{ Row = 9;  Values = @(\"48\u00f78=6, 0\", \"24\u00f72=12, 0\", \"99\u00f78=12, 3\", \"43\u00f74=10, 3\", \"67\u00f75=13, 2\") },\n    @{ Row = 13; Values = @(\"69\u00f75=13, 4\", \"37\u00f73=12, 1\", \"82\u00f78=10, 2\", \"60\u00f76=10, 0\", \"97\u00f74=24, 1\") },\n    @{ Row = 17; Values = @(\"51\u00f77=7, 2\", \"84\u00f74=21, 0\", \"18\u00f79=2, 0\", \"39\u00f79=4, 3\", \"69\u00f78=8, 5\") }\n)\n\nforeach ($rowRep in $rowReplacements) {\n    $rowIndex = $rowRep.Row\n    $values = $rowRep.Values\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
